$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E10").Value = 1

$ws.Range("E15").Value = 104
$ws.Range("F15").Value = 44
$ws.Range("H15").Value = 44

$ws.Range("E18").Value = 54
$ws.Range("F18").Value = 16
$ws.Range("H18").Value = 16

$ws.Range("E25").Value = 12

$ws.Range("E34").Value = 9

$ws.Range("E37").Value = 24
$ws.Range("F37").Value = 10
$ws.Range("H37").Value = 10

$ws.Range("E50").Value = 11

$ws.Range("E62").Value = 19

$ws.Range("E70").Value = 18
$ws.Range("F70").Value = 6
$ws.Range("H70").Value = 6

$ws.Range("E71").Value = 16

$ws.Range("E72").Value = 20
$ws.Range("F72").Value = 10
$ws.Range("H72").Value = 10

$ws.Range("E79").Value = 15

$ws.Range("E89").Value = 15
